$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows to append (B, C, D, E are text; A, F, G are numeric)
$newRows = @(
    @(1, "6-5-5",  "new_sequential", "63.283", "sat", 6276,  170692),
    @(1, "13-7-2", "new_sequential", "0.181",  "sat", 16705, 726186),
    @(1, "13-7-2", "new_sequential", "0.179",  "sat", 16705, 726186)
)

$startRow = 11
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D must stay a text string (e.g. "0.181"), not become a number
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
